# Regenerate the "K" (strikeouts) column (G) of save_data using the
# corrected per-start values instead of the old "Strike#" figures.
# This also corresponds to the regen of std/mean + s_vals noted in the
# commit message; only the stored K values themselves change here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(2,4,4,7,3,3,2,8,7,3,10,8,4,10,2,2,5,5,3,2,4,4,4,4,5,6,7,2,4,3,5,5,4,3,6,1)

$firstRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
